# Adds two new weekly survey periods (7.–13. 6. 2021 and 14.–20. 6. 2021) to
# both data sheets, and bumps the "aktualizace" (last-updated) date in the
# two footnote cells from 1. 6. 2021 to 28. 6. 2021.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": header row 1 spans A:BG -> new weeks go to columns BH/BI
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

$wsData.Cells.Item(1, 59).Copy()
$wsData.Cells.Item(1, 60).PasteSpecial(-4122)
$wsData.Cells.Item(1, 60).Value = "7.–13. 6. 2021"

$wsData.Cells.Item(1, 59).Copy()
$wsData.Cells.Item(1, 61).PasteSpecial(-4122)
$wsData.Cells.Item(1, 61).Value = "14.–20. 6. 2021"

$dataRows = @"
2,19,20
3,4,4
4,11.5,11.5
5,2.5,2.5
6,24.5,25.5
7,23.5,24
8,12.5,13
9,24,25
10,18,19
11,15.5,16
12,26,27
13,18.5,20
14,17,17.5
15,20.5,21.5
16,18,18.5
17,19,19.5
18,20,21.5
19,19.5,20.5
20,18.5,19
21,29,29
22,18.5,19.5
23,12.5,13.5
24,9.5,10.5
25,7.5,8
26,4.5,4.5
27,1.5,1.5
28,6.5,6.5
29,3.5,3.5
30,2,2.5
31,5,4.5
32,3,3.5
33,4,4
34,4,4
35,4,4
36,3.5,4
37,4.5,5
38,5,5
39,2.5,2.5
40,5.5,5.5
41,3,3.5
42,4,4
43,1,1.5
44,15.5,15.5
45,13,13
46,7.5,7.5
47,15,15
48,10.5,10.5
49,8,8
50,15.5,15.5
51,9.5,9.5
52,10.5,10.5
53,12.5,12.5
54,10,10
55,11.5,11.5
56,11.5,11.5
57,12,12.5
58,10,10
59,15,15
60,11.5,11.5
61,7.5,7.5
62,7.5,7.5
63,5,5
64,2.5,2.5
65,1,1
66,4.5,4.5
67,2,2
68,1,1
69,3.5,3.5
70,2,2
71,2.5,2.5
72,2.5,2.5
73,2.5,2.5
74,2.5,2.5
75,2.5,3
76,3,3
77,1.5,1.5
78,3,3.5
79,2,1.5
80,1.5,1.5
81,1,1
"@

foreach ($line in ($dataRows -split "`n")) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $p = $line -split ","
    $r = [int]$p[0]
    $wsData.Cells.Item($r, 60).Value = [double]$p[1]
    $wsData.Cells.Item($r, 61).Value = [double]$p[2]
}

# Footnote in A82: "... aktualizace 1. 6. 2021" -> "... aktualizace 28. 6. 2021"
$wsData.Cells.Item(82, 1).Value = "Život během pandemie, Kontakt s lidmi, průměr celkově a ve skupinách, aktualizace 28. 6. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": header row 1 spans A:BF -> new weeks go to columns BG/BH
# ---------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

$wsPocet.Cells.Item(1, 58).Copy()
$wsPocet.Cells.Item(1, 59).PasteSpecial(-4122)
$wsPocet.Cells.Item(1, 59).Value = "7.–13. 6. 2021"

$wsPocet.Cells.Item(1, 58).Copy()
$wsPocet.Cells.Item(1, 60).PasteSpecial(-4122)
$wsPocet.Cells.Item(1, 60).Value = "14.–20. 6. 2021"

$pocetRows = @"
2,1686,1689
3,394,394
4,607,602
5,685,693
6,528,522
7,687,699
8,471,468
9,339,338
10,390,386
11,957,965
12,843,847
13,843,842
14,870,879
15,389,386
16,201,197
17,226,227
18,722,719
19,81,84
20,84,84
21,58,54
"@

foreach ($line in ($pocetRows -split "`n")) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $p = $line -split ","
    $r = [int]$p[0]
    $wsPocet.Cells.Item($r, 59).Value = [double]$p[1]
    $wsPocet.Cells.Item($r, 60).Value = [double]$p[2]
}

# Footnote in A22: "... aktualizace 1. 6. 2021" -> "... aktualizace 28. 6. 2021"
$wsPocet.Cells.Item(22, 1).Value = "Život během pandemie, Kontakt s lidmi, velikost dotázaného souboru celkově a ve skupinách, aktualizace 28. 6. 2021"
